# Add a new column BB to the YoY forecast table, mirroring column BA for
# the historical rows (3-18) and the header date (row 1), while rows
# 19-21 (the "live"/rolling forecast rows) get newly computed values
# that differ from column BA.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: new forecast-vintage date header. Copy BA1 (value + style) then
# overwrite with the new date serial so formatting/border match the rest
# of the header row.
$ws.Range("BA1").Copy($ws.Range("BB1"))
$ws.Range("BB1").Value = 45986

# Rows 3-18: same YoY values as column BA (straight carry-forward).
$ws.Range("BA3").Copy($ws.Range("BB3"))
$ws.Range("BA4").Copy($ws.Range("BB4"))
$ws.Range("BA5").Copy($ws.Range("BB5"))
$ws.Range("BA6").Copy($ws.Range("BB6"))
$ws.Range("BA7").Copy($ws.Range("BB7"))
$ws.Range("BA8").Copy($ws.Range("BB8"))
$ws.Range("BA9").Copy($ws.Range("BB9"))
$ws.Range("BA10").Copy($ws.Range("BB10"))
$ws.Range("BA11").Copy($ws.Range("BB11"))
$ws.Range("BA12").Copy($ws.Range("BB12"))
$ws.Range("BA13").Copy($ws.Range("BB13"))
$ws.Range("BA14").Copy($ws.Range("BB14"))
$ws.Range("BA15").Copy($ws.Range("BB15"))
$ws.Range("BA16").Copy($ws.Range("BB16"))
$ws.Range("BA17").Copy($ws.Range("BB17"))
$ws.Range("BA18").Copy($ws.Range("BB18"))

# Rows 19-21: recomputed forecast values (differ from BA).
$ws.Range("BA19").Copy($ws.Range("BB19"))
$ws.Range("BB19").Value = 2.043309689777173

$ws.Range("BA20").Copy($ws.Range("BB20"))
$ws.Range("BB20").Value = 1.199077969291551

$ws.Range("BA21").Copy($ws.Range("BB21"))
$ws.Range("BB21").Value = 1.308808515504123

Write-Output "Added column BB (through row 21) to sheet $($ws.Name)"
